$wb = $excel.ActiveWorkbook

$wsInscrit = $wb.Worksheets.Item("Inscrit")
$wsTotal = $wb.Worksheets.Item("Total")

# --- Sheet "Inscrit" ---
# Row 2: remove the "x" previously in C2, move it to E2 (A2/B2 "t" stay as-is)
$wsInscrit.Range("C2").ClearContents()
$wsInscrit.Range("E2").Value = "x"

# Row 3 (new entry)
$wsInscrit.Range("A3").Value = "NAKHIL"
$wsInscrit.Range("B3").Value = "Amine"

# Row 4 (new entry)
$wsInscrit.Range("A4").Value = "l"
$wsInscrit.Range("B4").Value = "l"
$wsInscrit.Range("C4").Value = "x"

# --- Sheet "Total" ---
$wsTotal.Range("C2").Value = 1
$wsTotal.Range("E2").Value = 2
